$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.036.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.13%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.584.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.75%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("E5").Value = "  +0.01%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "301.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.04%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3761"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.16%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3583"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.53%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.62"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.46%  "

# Row 10
$ws.Range("E10").Value = "  +0.06%  "

# Row 11
$ws.Range("E11").Value = "  -4.14%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08038"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.07%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.10%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.463"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.19%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.307"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.21%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001222"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.23%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.589.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.26%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.87%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06807"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.25%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.28%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.09%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.434"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.11%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.35%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.040.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.18%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.360"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.09%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.782"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.33%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.27%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "147.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.05%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.206"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.07%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.63"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.67%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.354"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.30%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.535"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.17%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.757.67"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.80%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9412"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.67%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07356"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.52%  "

# Row 36
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02675"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.67%  "

# Row 37
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.32%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08754"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.17%  "

# Row 39
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.049"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.41%  "

# Row 40
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2473"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.31%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.333"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.22%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6874"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.56%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.99%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.37%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6403"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.94%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.993"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.78%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.246"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.44%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.25%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07878"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.29%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.192"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.04%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.202"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.62%  "
